# Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study
#
# The CasesTab "StatQuery" (column B, row 2 / "CasesTab") no longer returns
# the Cohort column - drop the trailing `co.cohort_description` RETURN line
# (and its now-dangling trailing comma on the previous line).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesTabQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.sex IN ['Unknown']
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $casesTabQuery

# Move the active selection from B4 to B2 (also scrolls the view back to
# the top, dropping the old topLeftCell="A4" anchor).
$ws.Range("B2").Select() | Out-Null
